# Trade #4 closed at 2026-02-17 15:13:38 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook to record the close of trade #4:
#   - Summary sheet roll-up metrics
#   - Strategy Status sheet's MarketMaking row
#   - A new trade row appended to both the "All Trades" and "MarketMaking" sheets

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.99   # Current Capital
$summary.Range("B4").Value = -0.01     # Total P&L $
$summary.Range("B5").Value = -0.05     # Total P&L %
$summary.Range("B6").Value = 4         # Total Trades
$summary.Range("B8").Value = 2         # Losing Trades
$summary.Range("B9").Value = 50        # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.98999999999999   # Capital
$status.Range("D4").Value = 4                   # Trades
$status.Range("E4").Value = -0.01               # P&L $
$status.Range("F4").Value = -0.01               # P&L %
$status.Range("G4").Value = 50                  # Win Rate %

# --- New trade row (#4) shared by "All Trades" and "MarketMaking" sheets ---
$tradeRow = @(
    4,
    "2026-02-17",
    "15:13:32",
    "MarketMaking",
    "UP",
    0.84,
    0.83,
    "CLOSED",
    -1.1905,
    -0.01,
    99.98999999999999,
    0,
    0,
    0.6,
    "Normal spread capture: 19600 bps",
    "early_exit",
    0.13
)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($col = 1; $col -le $tradeRow.Length; $col++) {
        $cell = $ws.Cells.Item(5, $col)
        $value = $tradeRow[$col - 1]

        if ($col -eq 2) {
            # Column B holds an ISO date string ("2026-02-17"). A bare
            # assignment gets auto-parsed into a date serial by Excel's
            # input parser, so force literal text (leading apostrophe),
            # then drop the resulting "quote prefix" cell format so the
            # cell ends up as plain text with no special styling.
            $cell.Value = "'" + $value
            $cell.ClearFormats()
        } else {
            $cell.Value = $value
        }
    }
}
